$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2: serial no bumped, task text simplified
Set-TextValue "B2" "2"
Set-TextValue "C2" "dfg"
Set-TextValue "D2" "dfg"

# Row 3
Set-TextValue "B3" "3"
Set-TextValue "C3" "fdg"
Set-TextValue "D3" "dfg"

# Row 4
Set-TextValue "B4" "4"
Set-TextValue "C4" "dfg"
Set-TextValue "D4" "dfg"

# New row 5 - copy the formatting of A4 (bordered/centered style) onto A5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").Value = 4

Set-TextValue "B5" "5"
Set-TextValue "C5" "dfg"
Set-TextValue "D5" "dfg"
